$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.215.85"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "1.602.48"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9985"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.09"
$ws.Range("E6").Value = "  +0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3779"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.97"
$ws.Range("E8").Value = "  +4.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3616"
$ws.Range("E9").Value = "  -0.81%  "

$ws.Range("E10").Value = "  +0.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08125"
$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9986"
$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.69"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.600"
$ws.Range("E14").Value = "  -0.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.405"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001252"
$ws.Range("E16").Value = "  -0.32%  "

$ws.Range("D17").Value = "1.603.45"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.41"
$ws.Range("E18").Value = "  +1.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06869"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.06"
$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.525"
$ws.Range("E21").Value = "  -0.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.97"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").Value = "23.219.66"
$ws.Range("E24").Value = "  +0.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.394"
$ws.Range("E25").Value = "  +2.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.000"
$ws.Range("E26").Value = "  +7.44%  "

$ws.Range("E27").Value = "  +0.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.84"
$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.230"
$ws.Range("E29").Value = "  -1.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.98"
$ws.Range("E30").Value = "  +1.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.402"
$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.806"
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("D33").Value = "1.778.05"
$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9814"
$ws.Range("E34").Value = "  +3.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07603"
$ws.Range("E35").Value = "  -1.16%  "

$ws.Range("E36").Value = "  +2.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02724"
$ws.Range("E37").Value = "  -1.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.160"
$ws.Range("E38").Value = "  -1.58%  "

$ws.Range("E39").Value = "  -1.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08797"
$ws.Range("E40").Value = "  -1.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7142"
$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.363"
$ws.Range("E42").Value = "  -2.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.43"
$ws.Range("E43").Value = "  -2.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.45"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6582"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.307"
$ws.Range("E46").Value = "  +0.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.017"
$ws.Range("E47").Value = "  +1.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.25"
$ws.Range("E48").Value = "  -0.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07961"
$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.222"
$ws.Range("E51").Value = "  +3.05%  "
